$d = $word.ActiveDocument

# The document contains two occurrences of the merge-field placeholder
# "{fechaEliminacion}" (the 2nd one is immediately followed by a period).
# Both get renamed to "{fechaEliminacionTexto}" -- but the rename is
# expressed in the OOXML as brand-new <w:r> runs appended right after the
# existing "fechaEliminacion"/"Eliminacion" text (not as a rewrite of the
# existing runs' text), so we insert new runs surgically instead of using
# a simple Find/Replace (which would merge everything into one run).

$search = $d.Content

# --- Occurrence #1: "Medellin, {fechaEliminacion}" ---------------------
$null = $search.Find.Execute("{fechaEliminacion}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insertAt1 = $d.Range($search.End - 1, $search.End - 1)
$insertAt1.InsertAfter("Texto")
# Round-trip a character property on just the inserted text so the engine
# materializes it as its own run (empty <w:rPr/>) instead of silently
# merging it into the preceding "fechaEliminacion" run.
$insertAt1.Bold = 1
$insertAt1.Bold = 0

# Keep searching after occurrence #1 so the next Find lands on #2.
$search.Collapse(0)

# --- Occurrence #2: "... desde {fechaEliminacion}." ---------------------
$null = $search.Find.Execute("{fechaEliminacion}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insertAtT = $d.Range($search.End - 1, $search.End - 1)
$insertAtT.InsertAfter("T")
$posExto = $insertAtT.End
$insertAtExto = $d.Range($posExto, $posExto)
$insertAtExto.InsertAfter("exto")

# Materialize both new runs as distinct <w:r> elements. Do the
# right-hand ("exto") run's round-trip FIRST and the left-hand ("T")
# run's round-trip LAST -- doing it in the opposite order lets the
# engine's run-consolidation pass re-merge the earlier split.
$extoRange = $d.Range($posExto, $posExto + 4)
$extoRange.Bold = 1
$extoRange.Bold = 0

$tRange = $d.Range($insertAtT.Start, $insertAtT.Start + 1)
$tRange.Bold = 1
$tRange.Bold = 0

Write-Output "done"
